$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 6895
$ws.Range("K3").Value = 7140
$ws.Range("D4").Value = 1979
$ws.Range("I4").Value = 1808
$ws.Range("K4").Value = 1476
$ws.Range("K5").Value = 509
$ws.Range("K6").Value = 7846
$ws.Range("D7").Value = 28169
$ws.Range("I7").Value = 26268
$ws.Range("K7").Value = 23866

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 438
$ws.Range("K3").Value = 475
$ws.Range("K6").Value = 518
$ws.Range("K7").Value = 1565

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 175
$ws.Range("K7").Value = 508

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 365
$ws.Range("K7").Value = 1024

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K3").Value = 137
$ws.Range("K7").Value = 393

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 233
$ws.Range("K3").Value = 265
$ws.Range("K7").Value = 807

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 182
$ws.Range("K4").Value = 21
$ws.Range("K6").Value = 201
$ws.Range("K7").Value = 552

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K6").Value = 172
$ws.Range("K7").Value = 729
$ws.Range("K8").Value = 1565
$ws.Range("K9").Value = 108
$ws.Range("K15").Value = 249
$ws.Range("K16").Value = 59
$ws.Range("K19").Value = 705
$ws.Range("K20").Value = 578
$ws.Range("K24").Value = 73
$ws.Range("K29").Value = 1300
$ws.Range("K33").Value = 1024
$ws.Range("K36").Value = 298
$ws.Range("K37").Value = 807
$ws.Range("K39").Value = 28
$ws.Range("K41").Value = 168
$ws.Range("K42").Value = 883
$ws.Range("K48").Value = 310
$ws.Range("K51").Value = 301
$ws.Range("K54").Value = 464
$ws.Range("D63").Value = 358
$ws.Range("I63").Value = 229
$ws.Range("K63").Value = 67
$ws.Range("K65").Value = 552
$ws.Range("K66").Value = 73
$ws.Range("K67").Value = 924
$ws.Range("K71").Value = 74
$ws.Range("K75").Value = 74
$ws.Range("K76").Value = 320
$ws.Range("K78").Value = 280
$ws.Range("K80").Value = 85
$ws.Range("K81").Value = 17
$ws.Range("K83").Value = 508
$ws.Range("K85").Value = 1100
$ws.Range("K86").Value = 149
$ws.Range("K90").Value = 227
$ws.Range("K91").Value = 282
$ws.Range("K93").Value = 90
$ws.Range("K94").Value = 319
$ws.Range("K95").Value = 393
$ws.Range("K96").Value = 257
$ws.Range("K98").Value = 121
$ws.Range("D101").Value = 28169
$ws.Range("I101").Value = 26268
$ws.Range("K101").Value = 23866

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 251
$ws.Range("K6").Value = 265
$ws.Range("K7").Value = 924

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K6").Value = 253
$ws.Range("K7").Value = 464

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 367
$ws.Range("K3").Value = 464
$ws.Range("K7").Value = 1300

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K3").Value = 73
$ws.Range("K7").Value = 310

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K6").Value = 234
$ws.Range("K7").Value = 705

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K6").Value = 162
$ws.Range("K7").Value = 320

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K3").Value = 46
$ws.Range("K6").Value = 45
$ws.Range("K7").Value = 172

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K6").Value = 68
$ws.Range("K7").Value = 168

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 240
$ws.Range("K6").Value = 330
$ws.Range("K7").Value = 883

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K6").Value = 94
$ws.Range("K7").Value = 280

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("K2").Value = 30
$ws.Range("K7").Value = 73

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K4").Value = 16
$ws.Range("K6").Value = 107
$ws.Range("K7").Value = 257

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K3").Value = 134
$ws.Range("K7").Value = 282

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 199
$ws.Range("K6").Value = 158
$ws.Range("K7").Value = 578

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K6").Value = 67
$ws.Range("K7").Value = 298

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 90

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K4").Value = 27
$ws.Range("K7").Value = 729

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K3").Value = 66
$ws.Range("K7").Value = 319

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K4").Value = 16
$ws.Range("K7").Value = 249

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K6").Value = 71
$ws.Range("K7").Value = 121

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("K5").Value = 18
$ws.Range("K6").Value = 28

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 73

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K2").Value = 35
$ws.Range("K7").Value = 108

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K2").Value = 26
$ws.Range("K3").Value = 23
$ws.Range("K4").Value = 63
$ws.Range("K5").Value = 2
$ws.Range("K7").Value = 149

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("K2").Value = 25
$ws.Range("K7").Value = 74

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 83
$ws.Range("K3").Value = 65
$ws.Range("K7").Value = 227

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K4").Value = 32
$ws.Range("K7").Value = 301

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 362
$ws.Range("K6").Value = 268
$ws.Range("K7").Value = 1100

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("K2").Value = 23
$ws.Range("K7").Value = 74

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 85

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("K6").Value = 34
$ws.Range("K7").Value = 59

$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 17
